# --- Append new scrape run: 2025-09-23 18:26 JST ---
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: newly scraped listing, inserted at the top
$ws.Range("A2").Value = "2025-09-23 18:26:38"
$ws.Range("B2").Value = "UbersuggestにてWEB上のデータを自動取得 Python マクロ RPA AIエージェント"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5399085"
$ws.Range("G2").Value = 490
$ws.Range("H2").Value = "🔥AI,Python"

# Rows 3-10: previous rows 2-9 shifted down one, timestamp refreshed
# Row 3
$ws.Range("A3").Value = "2025-09-23 18:26:38"
$ws.Range("B3").Value = "初回 「AIで笑顔を検出し、2秒クリップを無劣化で自動切り出すWindowsツール開発(予算10万円)」"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5398662"
$ws.Range("G3").Value = 413
$ws.Range("H3").Value = "🔥AI,Ai ◆ツール,開発"

# Row 4
$ws.Range("A4").Value = "2025-09-23 18:26:38"
$ws.Range("B4").Value = "【急募】出品・在庫管理ツール開発と保守対応者募集"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5398562"
$ws.Range("G4").Value = 163
$ws.Range("H4").Value = "◆ツール,開発 ◇管理"

# Row 5
$ws.Range("A5").Value = "2025-09-23 18:26:38"
$ws.Range("B5").Value = "海外仕入れ元サイト→ツールを動かす為のCSVファイルに週1で自動抽出の制作(自動/スクレイピング)"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5251319"
$ws.Range("G5").Value = 135
$ws.Range("H5").Value = "◆ツール,スクレイピング ◇サイト"

# Row 6
$ws.Range("A6").Value = "2025-09-23 18:26:38"
$ws.Range("B6").Value = "【急募】冠婚葬祭業公式サイトのPHP+MySQLバージョンアップ依頼"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5385491"
$ws.Range("G6").Value = 100
$ws.Range("H6").Value = "◇MySQL ○PHP"

# Row 7
$ws.Range("A7").Value = "2025-09-23 18:26:38"
$ws.Range("B7").Value = "EC多プラットフォーム展開在庫・価格連携ツールの作成"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5398432"
$ws.Range("G7").Value = 80
$ws.Range("H7").Value = "◆ツール"

# Row 8
$ws.Range("A8").Value = "2025-09-23 18:26:38"
$ws.Range("B8").Value = "運送会社の作業予定表 WEBシステムの修正開発についての相談"
$ws.Range("C8").Value = "システム開発"
$ws.Range("D8").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E8").Value = "期限情報なし"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5398772"
$ws.Range("G8").Value = 78
$ws.Range("H8").Value = "◆開発"

# Row 9
$ws.Range("A9").Value = "2025-09-23 18:26:38"
$ws.Range("B9").Value = "仮想通貨トレードの運用とコンサル【1名】のみ募集"
$ws.Range("C9").Value = "システム開発"
$ws.Range("D9").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E9").Value = "期限情報なし"
$ws.Range("F9").Value = "https://www.lancers.jp/work/detail/5393695"
$ws.Range("G9").Value = 55
$ws.Range("H9").Value = "◆コンサル"

# Row 10
$ws.Range("A10").Value = "2025-09-23 18:26:38"
$ws.Range("B10").Value = "【急募】災害時に備えた「ピジョンクラウド」でのシステムづくり、運用サポートの依頼"
$ws.Range("C10").Value = "システム開発"
$ws.Range("D10").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E10").Value = "期限情報なし"
$ws.Range("F10").Value = "https://www.lancers.jp/work/detail/5398657"
$ws.Range("G10").Value = 33

# Rows 11-12: unchanged listings, only the timestamp refreshes
$ws.Range("A11").Value = "2025-09-23 18:26:38"
$ws.Range("A12").Value = "2025-09-23 18:26:38"

# Rebuild the F-column hyperlinks against their (possibly new) target URLs.
$ws.Range("F2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5399085") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5398662") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5398562") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5251319") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5385491") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5398432") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5398772") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5393695") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5398657") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5398736") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.lancers.jp/work/detail/5395809") | Out-Null
$ws.Range("F2:F12").Style = "Hyperlink"
